$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing data row (251) down into the
# new rows (252:255) so the new date cells in column A keep the same
# style (date number format, centered, border) as the rest of the column.
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)

# Row 252: 2021-05-10
$ws.Range("A252").Value = 44326
$ws.Range("B252").Value = 0
$ws.Range("C252").Value = 0
$ws.Range("D252").Value = 0

# Row 253: 2021-05-11
$ws.Range("A253").Value = 44327
$ws.Range("B253").Value = 0
$ws.Range("C253").Value = 0
$ws.Range("D253").Value = 0

# Row 254: 2021-05-12
$ws.Range("A254").Value = 44328
$ws.Range("B254").Value = 1
$ws.Range("C254").Value = 1
$ws.Range("D254").Value = 37.46721618583739

# Row 255: 2021-05-13
$ws.Range("A255").Value = 44329
$ws.Range("B255").Value = 0
$ws.Range("C255").Value = 1
$ws.Range("D255").Value = 37.46721618583739
